$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text in B1: "Value" -> "Value (g)"
$ws.Range("B1").Value = "Value (g)"

# Move/collapse the selection from the whole-column A1:A1048576 down to D9
$ws.Range("D9").Select()
